$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.462.81"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").Value = "3.692.09"
$ws.Range("E3").Value = "  +0.75%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "685.82"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +1.29%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.11"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.01%  "

$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("D13").Value = "4.314.11"
$ws.Range("E13").Value = "  +0.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.55"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.27%  "

$ws.Range("D15").Value = "3.698.86"
$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("D16").Value = "69.453.19"
$ws.Range("E16").Value = "  +0.70%  "

$ws.Range("E17").Value = "  +1.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.88"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.43"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.09%  "

$ws.Range("E21").Value = "  +3.14%  "

$ws.Range("E22").Value = "  -1.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.77"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.25%  "

$ws.Range("D24").Value = "3.836.72"
$ws.Range("E24").Value = "  +0.53%  "

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("E26").Value = "  -0.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.04"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.28"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.35%  "

$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.75"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.01"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.59"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.97"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("D35").Value = "3.664.45"
$ws.Range("E35").Value = "  +0.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.159"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.21"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.18"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.17%  "

$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.23"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0902"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.17%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.944"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "165.88"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.57"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000284"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.94%  "

$ws.Range("E47").Value = "  +8.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.72%  "

$ws.Range("E49").Value = "  +1.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.92"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.78"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.23%  "
